$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1, matching the style used by the other
# header cells (e.g. G1 - bold font, border, centered alignment).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add corresponding numeric data value in H2 (no special style, like B2:G2).
$ws.Range("H2").Value = 0
